$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column E, matching style of existing header row (bold + centered)
$ws.Range("E1").Value = "Colocação"
$ws.Range("E1").Font.Bold = $true
$ws.Range("E1").HorizontalAlignment = -4108  # xlCenter

# Ranking values for rows 2-7 (states), rows 8-9 (Nordeste/Brasil totals) left blank
$ws.Range("E2").Value = "1º"
$ws.Range("E3").Value = "2º"
$ws.Range("E4").Value = "3º"
$ws.Range("E5").Value = "4º"
$ws.Range("E6").Value = "5º"
$ws.Range("E7").Value = "6º"
